$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("n9",  "n9_IMG_3177.jpeg",  "'True", "no_meltpatch", "negative"),
    @("n10", "n10_IMG_3175.jpeg", "'True", "no_meltpatch", "negative"),
    @("n11", "n11_IMG_3176.jpeg", "'True", "no_meltpatch", "negative"),
    @("n12", "n12_IMG_3178.jpeg", "'True", "no_meltpatch", "negative")
)

$startRow = 10
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $rowData[$col - 1]
        $cell.ClearFormats()
    }
}
